$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 18): A18 = 21, B18 = new TODO item text
$ws.Range("A18").Value = 21
$ws.Range("B18").Value = "Constraints (point-edge, face-face, perp, parallel, etc)"

# Move the active selection to B19, matching the post-edit workbook state
$ws.Range("B19").Select()
